$wb = $excel.ActiveWorkbook

# --- Metadata sheet: bump the generation Date ---
$meta = $wb.Worksheets.Item("Metadata")
$meta.Cells.Item(8, 2).Value = "2025-05-03T11:17:55+00:00"

# --- Elements sheet: rename Base Path entries PatientId.1/.4/.5 -> PatientId.CX1/.CX4/.CX5 ---
$elements = $wb.Worksheets.Item("Elements")

# Row 3: ID / Path / Base Path columns (A, B, AF) all mirror the element id
$elements.Cells.Item(3, 1).Value = "PatientId.CX1"
$elements.Cells.Item(3, 2).Value = "PatientId.CX1"
$elements.Cells.Item(3, 32).Value = "PatientId.CX1"

# Row 4
$elements.Cells.Item(4, 1).Value = "PatientId.CX4"
$elements.Cells.Item(4, 2).Value = "PatientId.CX4"
$elements.Cells.Item(4, 32).Value = "PatientId.CX4"

# Row 5
$elements.Cells.Item(5, 1).Value = "PatientId.CX5"
$elements.Cells.Item(5, 2).Value = "PatientId.CX5"
$elements.Cells.Item(5, 32).Value = "PatientId.CX5"

# The wider id text (".CXn" vs ".n") makes columns A, B and AF grow to fit
# (bestFit columns), so re-apply the best-fit width on them.
$elements.Columns.Item(1).ColumnWidth = 10.666666666666666
$elements.Columns.Item(2).ColumnWidth = 10.666666666666666
$elements.Columns.Item(32).ColumnWidth = 10.666666666666666

# Keep the existing frozen header/id panes on "Elements" (split after column
# B / row 1) intact, then restore "Metadata" as the active sheet/tab.
$elements.Activate() | Out-Null
$elements.Range("C2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true
$elements.Range("A2").Select() | Out-Null
$meta.Activate() | Out-Null
